$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to reflect the latest scrape
$ws.Range("D2").Value = "54.853.58"
$ws.Range("E2").Value = "  -2.72%  "

$ws.Range("D3").Value = "2.346.22"
$ws.Range("E3").Value = "  -5.72%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "470.67"
$ws.Range("E5").Value = "  -3.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.82"
$ws.Range("E6").Value = "  -2.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.501"
$ws.Range("E8").Value = "  -2.71%  "

$ws.Range("D9").Value = "2.344.13"
$ws.Range("E9").Value = "  -6.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0956"
$ws.Range("E10").Value = "  -2.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.38"
$ws.Range("E11").Value = "  -7.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.317"
$ws.Range("E12").Value = "  -4.89%  "

$ws.Range("E13").Value = "  +0.59%  "

$ws.Range("D14").Value = "2.757.84"
$ws.Range("E14").Value = "  -5.51%  "

$ws.Range("D15").Value = "54.940.50"
$ws.Range("E15").Value = "  -2.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.92"
$ws.Range("E16").Value = "  -6.35%  "

$ws.Range("E17").Value = "  -5.27%  "

$ws.Range("D18").Value = "2.356.61"
$ws.Range("E18").Value = "  -5.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.51"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "310.52"
$ws.Range("E20").Value = "  -2.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.53"
$ws.Range("E21").Value = "  -5.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.58"
$ws.Range("E23").Value = "  -4.38%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "55.78"
$ws.Range("E24").Value = "  -5.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.389"
$ws.Range("E26").Value = "  -5.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.153"
$ws.Range("E27").Value = "  -6.11%  "

$ws.Range("D28").Value = "2.451.83"
$ws.Range("E28").Value = "  -5.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.12"
$ws.Range("E29").Value = "  -7.08%  "

$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("D31").Value = "0.0₃0746"
$ws.Range("E31").Value = "  -5.61%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "146.88"
$ws.Range("E32").Value = "  -1.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.92"
$ws.Range("E33").Value = "  -1.82%  "

$ws.Range("E34").Value = "  -3.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.00"
$ws.Range("E35").Value = "  -4.58%  "

$ws.Range("E36").Value = "  -6.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.53"
$ws.Range("E37").Value = "  -5.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.815"
$ws.Range("E38").Value = "  -5.33%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.36"
$ws.Range("E39").Value = "  -2.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.33"
$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  -5.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0941"
$ws.Range("E43").Value = "  +2.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.574"
$ws.Range("E44").Value = "  -6.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0523"
$ws.Range("E45").Value = "  -7.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.16"
$ws.Range("E46").Value = "  -0.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "253.28"
$ws.Range("E47").Value = "  -2.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0219"
$ws.Range("E48").Value = "  -4.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.39"
$ws.Range("E49").Value = "  -8.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.65"
$ws.Range("E50").Value = "  -5.83%  "

$ws.Range("D51").Value = "1.773.88"
$ws.Range("E51").Value = "  -6.47%  "
